# Generate Report for Handoff
# Update the handoff/handback timestamps and set the "ht" (handoff type)
# priority value on the rows that were processed during this report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows whose latest handoff timestamp / priority need refreshing.
$rows = @(7, 8, 11, 12, 13, 14)

# Overview sheet: column G holds "Latest Handoff Datetime" -> bump from
# 2016-08-25 16:22:26 to 2016-08-25 16:22:41
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-25 16:22:41"
}

# zh-cn sheet: column H holds "Latest Handoff Datetime" -> bump from
# 2016-08-25 16:22:21 to 2016-08-25 16:22:37, and column E (Priority) is
# now set to "ht" for these rows.
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-25 16:22:37"
}

# de-de sheet: column H holds "Latest Handoff Datetime" -> bump from
# 2016-08-25 16:22:26 to 2016-08-25 16:22:41, and column E (Priority) is
# now set to "ht" for these rows.
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-25 16:22:41"
}
